# Applies the cryptos list refresh: updated prices/volumes for existing rows
# plus a reordering of four coin rows (44-47) to reflect the new ranking.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.322.61"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "3.775.61"
$ws.Range("E3").Value = "  -0.24%  "
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.84"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.19"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.49%  "
$ws.Range("D7").Value = "3.775.82"
$ws.Range("E7").Value = "  -0.33%  "
$ws.Range("E8").Value = "  -0.03%  "
$ws.Range("E9").Value = "  -1.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.161"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.51"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.449"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000265"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -3.53%  "
$ws.Range("E14").Value = "  -1.07%  "
$ws.Range("D15").Value = "4.406.32"
$ws.Range("E15").Value = "  -0.20%  "
$ws.Range("D16").Value = "3.771.40"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "68.322.01"
$ws.Range("E17").Value = "  +0.87%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "18.19"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.13%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.06"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.98%  "
$ws.Range("E20").Value = "  -0.40%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.05"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.45%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "468.64"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.703"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.14%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +1.33%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000144"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.50%  "
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.20"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.86%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("D30").Value = "3.921.69"
$ws.Range("E30").Value = "  -0.32%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.80"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.01%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.41"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.02%  "
$ws.Range("E33").Value = "  -1.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.09"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.26%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.998"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Value = "3.730.01"
$ws.Range("E37").Value = "  -0.46%  "
$ws.Range("E38").Value = "  -3.19%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.42"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -10.65%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.139"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.82%  "
$ws.Range("E41").Value = "  -0.61%  "
$ws.Range("E42").Value = "  -1.62%  "
$ws.Range("E43").Value = "  -0.03%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.309"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.56%  "
$ws.Range("B45").Value = "USDe"
$ws.Range("C45").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.00"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("B46").Value = "Arweave"
$ws.Range("C46").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "43.95"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +12.60%  "
$ws.Range("B47").Value = "Stacks"
$ws.Range("C47").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.96"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.31%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.59"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.76%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "407.30"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.44%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "45.52"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "145.50"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.41%  "

Write-Host "Applied 95 cell updates"
